$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '69.532.18'
$ws.Range('E2').Value = '  +2.13%  '

Set-TextCell 'D3' '3.350.47'
$ws.Range('E3').Value = '  +3.15%  '

Set-TextCell 'D5' '192.86'
$ws.Range('E5').Value = '  +4.64%  '

Set-TextCell 'D6' '596.48'
$ws.Range('E6').Value = '  +2.45%  '

$ws.Range('E7').Value = '  +0.02%  '

Set-TextCell 'D8' '0.608'
$ws.Range('E8').Value = '  +1.12%  '

$ws.Range('E9').Value = '  +3.33%  '

Set-TextCell 'D10' '6.72'
$ws.Range('E10').Value = '  +1.62%  '

Set-TextCell 'D11' '0.426'
$ws.Range('E11').Value = '  +2.25%  '

Set-TextCell 'D12' '3.928.27'
$ws.Range('E12').Value = '  +3.06%  '

$ws.Range('E13').Value = '  +0.93%  '

Set-TextCell 'D14' '28.39'
$ws.Range('E14').Value = '  +1.77%  '

Set-TextCell 'D15' '69.571.40'
$ws.Range('E15').Value = '  +2.15%  '

$ws.Range('E16').Value = '  +1.22%  '

Set-TextCell 'D17' '3.326.55'
$ws.Range('E17').Value = '  +1.53%  '

Set-TextCell 'D18' '5.86'
$ws.Range('E18').Value = '  +0.92%  '

Set-TextCell 'D19' '13.77'
$ws.Range('E19').Value = '  +2.12%  '

Set-TextCell 'D20' '427.55'
$ws.Range('E20').Value = '  +8.13%  '

Set-TextCell 'D21' '7.77'
$ws.Range('E21').Value = '  +2.12%  '

Set-TextCell 'D22' '72.91'
$ws.Range('E22').Value = '  +2.18%  '

$ws.Range('E23').Value = '  +0.23%  '

Set-TextCell 'D24' '0.521'
$ws.Range('E24').Value = '  +0.84%  '

Set-TextCell 'D25' '0.0000122'
$ws.Range('E25').Value = '  +2.99%  '

Set-TextCell 'D26' '0.192'
$ws.Range('E26').Value = '  +3.10%  '

Set-TextCell 'D27' '9.65'
$ws.Range('E27').Value = '  -0.04%  '

Set-TextCell 'D28' '1.01'
$ws.Range('E28').Value = '  +0.58%  '

Set-TextCell 'D29' '2.02'
$ws.Range('E29').Value = '  +2.34%  '

$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell 'D30' '5.70'
$ws.Range('E30').Value = '  +1.57%  '

$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextCell 'D31' '23.17'
$ws.Range('E31').Value = '  +1.60%  '

Set-TextCell 'D32' '1.29'
$ws.Range('E32').Value = '  +1.94%  '

Set-TextCell 'D33' '7.08'
$ws.Range('E33').Value = '  +0.75%  '

$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell 'D34' '0.999'
$ws.Range('E34').Value = '  -0.01%  '

$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D35' '164.19'
$ws.Range('E35').Value = '  +1.52%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell 'D36' '1.51'
$ws.Range('E36').Value = '  +0.99%  '

$ws.Range('B37').Value = 'Stacks'
$ws.Range('C37').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell 'D37' '1.94'
$ws.Range('E37').Value = '  +1.49%  '

$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D38' '27.35'
$ws.Range('E38').Value = '  +3.18%  '

$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D39' '4.62'
$ws.Range('E39').Value = '  +0.31%  '

$ws.Range('B40').Value = 'Mantle'
$ws.Range('C40').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 'D40' '0.811'
$ws.Range('E40').Value = '  -0.53%  '

$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D41' '6.50'
$ws.Range('E41').Value = '  -0.04%  '

$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell 'D42' '2.744.94'
$ws.Range('E42').Value = '  +5.28%  '

$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell 'D43' '2.53'
$ws.Range('E43').Value = '  +2.08%  '

Set-TextCell 'D44' '25.66'
$ws.Range('E44').Value = '  +2.58%  '

$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D45' '0.0690'
$ws.Range('E45').Value = '  +0.80%  '

Set-TextCell 'D46' '41.12'
$ws.Range('E46').Value = '  -0.25%  '

$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell 'D47' '344.94'
$ws.Range('E47').Value = '  +2.79%  '

$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 'D48' '0.0283'
$ws.Range('E48').Value = '  +1.17%  '

$ws.Range('B49').Value = 'Arweave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextCell 'D49' '32.74'
$ws.Range('E49').Value = '  +5.33%  '

$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextCell 'D50' '1.01'
$ws.Range('E50').Value = '  +3.48%  '

$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D51' '6.33'
$ws.Range('E51').Value = '  -0.10%  '

